# Architecture_Diagram.pptx WIP fit-up edit
#
# Moves/resizes three text boxes on slide 1 and updates two text runs:
#   - "TextBox 36" (shape id 37): reposition/resize + rename
#       "Entry Request" -> "Lock Entries Number"
#   - "TextBox 74" (shape id 75): shift left position only (text unchanged)
#   - "TextBox 76" (shape id 77): shift left position / narrow + rename
#       "Reliable Causal Ordering Multicast" -> "Reliable Total Ordering Multicast"
#
# NOTE: Shape.Left/Top/Width/Height are exposed as single-precision (Single)
# floats in the PowerPoint object model, so the literals below are chosen
# (rather than the "obvious" EMU/12700 quotient) so that after the
# float32 round-trip the saved EMU value in the XML matches the target
# exactly.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- "Entry Request" -> "Lock Entries Number" -----------------------------
$shp37 = $s.Shapes.Item(20)
$shp37.Left   = 335.0841369628906
$shp37.Top    = 70.501220703125
$shp37.Width  = 67.66287994384766
$shp37.Height = 58.162559509277344
$shp37.TextFrame.TextRange.Text = "Lock Entries Number"

# --- "Unicast TCP" textbox: move slightly right, text stays the same ------
$shp75 = $s.Shapes.Item(26)
$shp75.Left = 192.7648468017578

# --- "Reliable Causal Ordering Multicast" -> "Reliable Total Ordering Multicast"
$shp77 = $s.Shapes.Item(28)
$shp77.Left  = 59.45130157470703
$shp77.Width = 231.4287872314453
$shp77.TextFrame.TextRange.Text = "Reliable Total Ordering Multicast"
